$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 15:35"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1708265
$ws.Range("C4").Value = 2039
$ws.Range("E4").Value = 1143691
$ws.Range("G4").Value = 41
$ws.Range("H4").Value = 99846

# India (row 13)
$ws.Range("D13").Value = 61312
$ws.Range("E13").Value = 80998
$ws.Range("G13").Value = 16
$ws.Range("H13").Value = 4188

# Arabia Saudita (row 18)
$ws.Range("B18").Value = 76726
$ws.Range("C18").Value = 1931
$ws.Range("D18").Value = 48450
$ws.Range("E18").Value = 27865
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 411

# Suiza (row 32)
$ws.Range("E32").Value = 646
$ws.Range("G32").Value = 2
$ws.Range("H32").Value = 1915

# Argentina (row 47)
$ws.Range("D47").Value = 4167
$ws.Range("E47").Value = 7990
$ws.Range("G47").Value = 4
$ws.Range("H47").Value = 471

# Noruega (row 57)
$ws.Range("B57").Value = 8374
$ws.Range("C57").Value = 10
$ws.Range("E57").Value = 412

# Cuba (row 91)
$ws.Range("B91").Value = 1963
$ws.Range("C91").Value = 16
$ws.Range("D91").Value = 1709
$ws.Range("E91").Value = 172

# Islandia (row 93)
$ws.Range("D93").Value = 1792
$ws.Range("E93").Value = 2

# Kenia (row 102)
$ws.Range("B102").Value = 1348
$ws.Range("C102").Value = 62
$ws.Range("D102").Value = 405
$ws.Range("E102").Value = 891

# Sri Lanka (row 103)
$ws.Range("B103").Value = 1278
$ws.Range("C103").Value = 96
$ws.Range("E103").Value = 556

# Estado de Palestina (row 139)
$ws.Range("B139").Value = 426
$ws.Range("C139").Value = 3
$ws.Range("D139").Value = 365
$ws.Range("E139").Value = 58

$wb.Save()
